# Update computed profit/price figures in the Zeromus_Profits sheets
# (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) per the latest scheduled data refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H92").Value = 4630200
$ws.Range("I92").Value = 6173221
$ws.Range("J92").Value = 1137.1111
$ws.Range("K92").Value = 6173221
$ws.Range("L92").Value = 1137.1111
$ws.Range("M92").Value = -6171973
$ws.Range("N92").Value = -3633.1111

$ws.Range("H135").Value = 985.30554
$ws.Range("I135").Value = 615.59375
$ws.Range("K135").Value = 5540.34375
$ws.Range("M135").Value = -3005.34375

$ws.Range("H138").Value = 1742.4149
$ws.Range("I138").Value = 539.25
$ws.Range("J138").Value = 2997.8914
$ws.Range("K138").Value = 1617.75
$ws.Range("L138").Value = 8993.674199999999
$ws.Range("M138").Value = 3522.25
$ws.Range("N138").Value = -19273.6742

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9452.639999999999
$ws.Range("I32").Value = 2947.8057
$ws.Range("K32").Value = 2947.8057
$ws.Range("M32").Value = -2660.8057

$ws.Range("H110").Value = 16770.25
$ws.Range("I110").Value = 111111
$ws.Range("J110").Value = 3293
$ws.Range("K110").Value = 111111
$ws.Range("L110").Value = 3293
$ws.Range("M110").Value = -109066
$ws.Range("N110").Value = -7383

$ws.Range("H132").Value = 1965.4531
$ws.Range("I132").Value = 1320.7073
$ws.Range("J132").Value = 3114.7827
$ws.Range("K132").Value = 3962.1219
$ws.Range("L132").Value = 9344.348100000001
$ws.Range("M132").Value = -1432.1219
$ws.Range("N132").Value = -14404.3481

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1446.6904
$ws.Range("I134").Value = 1162.5555
$ws.Range("J134").Value = 3151.5
$ws.Range("K134").Value = 3487.6665
$ws.Range("L134").Value = 9454.5
$ws.Range("M134").Value = -952.6664999999998
$ws.Range("N134").Value = -14524.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H20").Value = 33299.832
$ws.Range("J20").Value = 33299.832
$ws.Range("L20").Value = 33299.832
$ws.Range("N20").Value = -33771.832

$ws.Range("H21").Value = 1974.75
$ws.Range("I21").Value = 1399
$ws.Range("J21").Value = 2166.6667
$ws.Range("K21").Value = 1399
$ws.Range("L21").Value = 2166.6667
$ws.Range("M21").Value = -1164
$ws.Range("N21").Value = -2636.6667

$ws.Range("H23").Value = 17770
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 17770
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = 17770
$ws.Range("M23").ClearContents()
$ws.Range("N23").Value = -18250

$ws.Range("H27").Value = 17770
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 17770
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 17770
$ws.Range("M27").ClearContents()
$ws.Range("N27").Value = -18154

$ws.Range("H29").Value = 29800
$ws.Range("J29").Value = 29800
$ws.Range("L29").Value = 29800
$ws.Range("N29").Value = -30386

$ws.Range("H30").Value = 33299.832
$ws.Range("J30").Value = 33299.832
$ws.Range("L30").Value = 33299.832
$ws.Range("N30").Value = -33481.832

$ws.Range("H31").Value = 3142705.2
$ws.Range("I31").Value = 7179702.5
$ws.Range("J31").Value = 2818.6667
$ws.Range("K31").Value = 7179702.5
$ws.Range("L31").Value = 2818.6667
$ws.Range("M31").Value = -7179407.5
$ws.Range("N31").Value = -3408.6667

$ws.Range("H34").Value = 3142705.2
$ws.Range("I34").Value = 7179702.5
$ws.Range("J34").Value = 2818.6667
$ws.Range("K34").Value = 7179702.5
$ws.Range("L34").Value = 2818.6667
$ws.Range("M34").Value = -7179500.5
$ws.Range("N34").Value = -3222.6667

$ws.Range("H123").Value = 26777.777
$ws.Range("J123").Value = 26777.777
$ws.Range("L123").Value = 26777.777
$ws.Range("N123").Value = -36577.777

$ws.Range("H124").Value = 10465.2
$ws.Range("I124").Value = 2500
$ws.Range("J124").Value = 15775.333
$ws.Range("K124").Value = 2500
$ws.Range("L124").Value = 15775.333
$ws.Range("M124").Value = -45
$ws.Range("N124").Value = -20685.333

$ws.Range("H125").Value = 13333
$ws.Range("J125").Value = 13333
$ws.Range("L125").Value = 13333
$ws.Range("N125").Value = -18253

$ws.Range("H128").Value = 33299.832
$ws.Range("J128").Value = 33299.832
$ws.Range("L128").Value = 33299.832
$ws.Range("N128").Value = -43259.832

$ws.Range("H129").Value = 21097.445
$ws.Range("J129").Value = 21097.445
$ws.Range("L129").Value = 21097.445
$ws.Range("N129").Value = -31097.445

$ws.Range("H130").Value = 40195
$ws.Range("J130").Value = 40195
$ws.Range("L130").Value = 40195
$ws.Range("N130").Value = -50235

$ws.Range("H131").Value = 29000
$ws.Range("J131").Value = 29000
$ws.Range("L131").Value = 29000
$ws.Range("N131").Value = -39080

$ws.Range("H132").Value = 1751.68
$ws.Range("I132").Value = 1264.7
$ws.Range("J132").Value = 3699.6
$ws.Range("K132").Value = 3794.1
$ws.Range("L132").Value = 11098.8
$ws.Range("M132").Value = -1264.1
$ws.Range("N132").Value = -16158.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1066.7142
$ws.Range("J131").Value = 1154.9762
$ws.Range("L131").Value = 3464.9286
$ws.Range("N131").Value = -13544.9286

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 1752.2903
$ws.Range("I126").Value = 1650.2916
$ws.Range("J126").Value = 2102
$ws.Range("K126").Value = 4950.8748
$ws.Range("L126").Value = 6306
$ws.Range("M126").Value = -2480.8748
$ws.Range("N126").Value = -11246

$ws.Range("H132").Value = 1693.6923
$ws.Range("I132").Value = 1297.7693
$ws.Range("J132").Value = 2485.5386
$ws.Range("K132").Value = 3893.3079
$ws.Range("L132").Value = 7456.6158
$ws.Range("M132").Value = -1363.3079
$ws.Range("N132").Value = -12516.6158

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 8933132
$ws.Range("I132").Value = 18389222
$ws.Range("J132").Value = 2380.861
$ws.Range("K132").Value = 55167666
$ws.Range("L132").Value = 7142.583
$ws.Range("M132").Value = -55165136
$ws.Range("N132").Value = -12202.583

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 18988.889
$ws.Range("J64").Value = 18988.889
$ws.Range("L64").Value = 18988.889
$ws.Range("N64").Value = -19484.889

$ws.Range("H67").Value = 18988.889
$ws.Range("J67").Value = 18988.889
$ws.Range("L67").Value = 18988.889
$ws.Range("N67").Value = -20704.889

$ws.Range("H122").Value = 1834.6666
$ws.Range("I122").Value = 1252
$ws.Range("J122").Value = 3000
$ws.Range("K122").Value = 3756
$ws.Range("L122").Value = 9000
$ws.Range("M122").Value = -1306
$ws.Range("N122").Value = -13900

$ws.Range("H123").Value = 46320.12
$ws.Range("J123").Value = 46320.12
$ws.Range("L123").Value = 46320.12
$ws.Range("N123").Value = -56120.12

$ws.Range("H132").Value = 1437.907
$ws.Range("I132").Value = 990.44446
$ws.Range("J132").Value = 2193
$ws.Range("K132").Value = 2971.33338
$ws.Range("L132").Value = 6579
$ws.Range("M132").Value = -441.33338
$ws.Range("N132").Value = -11639
